# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
#
# The existing sheet runs from column A to AC (29 columns); this appends
# three new columns - AD (Wins), AE (Losses), AF (Ties) - with a header
# row styled like the rest of the header, and a constant season record
# (75-87-0) for every player row (2 through 45).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Clone the formatting of an existing header cell (bold font, borders,
# centered/top alignment) onto the three new header cells, then set
# their text.
$ws.Range("A1").Copy($ws.Range("AD1"))
$ws.Range("A1").Copy($ws.Range("AE1"))
$ws.Range("A1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (rows 2-45) --------------------------------------------
# Every player row gets the team's season record: 75 wins, 87 losses,
# 0 ties.
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 75
    $ws.Cells.Item($row, 31).Value = 87
    $ws.Cells.Item($row, 32).Value = 0
}

Write-Host "Added Wins/Losses/Ties columns (AD:AF) for rows 1-45"
